$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-24 Friday", 2) | Out-Null
$d.Content.Find.Execute("85÷5=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=25, 0", 2) | Out-Null
$d.Content.Find.Execute("24÷3=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "25÷5=5, 0", 2) | Out-Null
$d.Content.Find.Execute("38÷2=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "43÷3=14, 1", 2) | Out-Null
$d.Content.Find.Execute("61÷3=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷2=5, 0", 2) | Out-Null
$d.Content.Find.Execute("94÷3=31, 1", $true, $false, $false, $false, $false, $true, 1, $false, "33÷6=5, 3", 2) | Out-Null
$d.Content.Find.Execute("12÷6=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=9, 4", 2) | Out-Null
$d.Content.Find.Execute("98÷7=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "52÷7=7, 3", 2) | Out-Null
$d.Content.Find.Execute("47÷3=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "11÷8=1, 3", 2) | Out-Null
$d.Content.Find.Execute("32÷7=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "80÷6=13, 2", 2) | Out-Null
$d.Content.Find.Execute("81÷4=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "20÷6=3, 2", 2) | Out-Null
$d.Content.Find.Execute("36÷4=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷6=11, 5", 2) | Out-Null
$d.Content.Find.Execute("42÷8=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "30÷8=3, 6", 2) | Out-Null
$d.Content.Find.Execute("83÷3=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "35÷3=11, 2", 2) | Out-Null
$d.Content.Find.Execute("66÷7=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "15÷4=3, 3", 2) | Out-Null
$d.Content.Find.Execute("10÷4=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=2, 5", 2) | Out-Null
$d.Content.Find.Execute("50÷7=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷2=39, 0", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "69÷5=13, 4", 2) | Out-Null
$d.Content.Find.Execute("67÷7=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=14, 3", 2) | Out-Null
$d.Content.Find.Execute("26÷9=2, 8", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=14, 1", 2) | Out-Null
$d.Content.Find.Execute("24÷8=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "35÷9=3, 8", 2) | Out-Null
$d.Content.Find.Execute("68÷8=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "62÷2=31, 0", 2) | Out-Null
$d.Content.Find.Execute("42÷7=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷2=44, 1", 2) | Out-Null
$d.Content.Find.Execute("19÷3=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=14, 2", 2) | Out-Null
$d.Content.Find.Execute("92÷7=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=7, 4", 2) | Out-Null
$d.Content.Find.Execute("91÷4=22, 3", $true, $false, $false, $false, $false, $true, 1, $false, "25÷9=2, 7", 2) | Out-Null
